# Apply odds updates to Sheet1 of the "Jogos da Semana FlashScore" workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2 ---
$ws.Range("N2").Value = 8

# --- Row 14 ---
$ws.Range("G14").Value = 2.57
$ws.Range("H14").Value = 3.2
$ws.Range("I14").Value = 2.5
$ws.Range("J14").Value = 3.2
$ws.Range("L14").Value = 3.1

$ws.Range("S14").Value = 1.4
$ws.Range("T14").Value = 2.72

$ws.Range("W14").Value = 8.25
$ws.Range("X14").Value = 13
$ws.Range("Y14").Value = 9.75
$ws.Range("Z14").Value = 29
$ws.Range("AA14").Value = 22

$ws.Range("AD14").Value = 6.3
$ws.Range("AE14").Value = 14

$ws.Range("AH14").Value = 8
$ws.Range("AI14").Value = 12.5
$ws.Range("AJ14").Value = 9.75
$ws.Range("AK14").Value = 27
$ws.Range("AL14").Value = 21

$ws.Range("AN14").Value = 4.55
$ws.Range("AO14").Value = 14
$ws.Range("AP14").Value = 22
$ws.Range("AQ14").Value = 60
$ws.Range("AR14").Value = 100
$ws.Range("AS14").Value = 300
$ws.Range("AT14").Value = 2.72
$ws.Range("AU14").Value = 7
$ws.Range("AV14").Value = 60
$ws.Range("AW14").Value = 4.45
$ws.Range("AX14").Value = 13.5
$ws.Range("AY14").Value = 21
$ws.Range("AZ14").Value = 55
$ws.Range("BA14").Value = 90
$ws.Range("BB14").Value = 250
